$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 5 query cells (B2=Program, C2=Count, B3=Project, B4=Grant, B5=Publication) ---
$countQuery = @'
SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Bladder Cancer%'
'@
$programQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE      prg.cancer_type LIKE '%Bladder Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@
$projectQuery = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Bladder Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@
$grantQuery = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.project_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Bladder Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@
$publicationQuery = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Bladder Cancer%'
ORDER BY 
    pub.pmid ASC
LIMIT 100;
'@

$ws.Range("C2").Value = $countQuery
$ws.Range("B2").Value = $programQuery
$ws.Range("B3").Value = $projectQuery
$ws.Range("B4").Value = $grantQuery
$ws.Range("B5").Value = $publicationQuery

# Ensure wrap-text formatting (font size 12, wrap) stays applied on B4/B5 matching the other query cells
$ws.Range("B4").WrapText = $true
$ws.Range("B4").Font.Size = 12
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Font.Size = 12
$ws.Range("C2").WrapText = $true
$ws.Range("C2").Font.Size = 12

# --- View state: scroll back to top-left A1 and move the active selection to C3 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C3").Select()
